$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.631.07"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "2.397.21"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.89"
$ws.Range("E5").Value = "  -1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.10"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("E8").Value = "  +1.51%  "

$ws.Range("D9").Value = "2.401.23"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("E10").Value = "  +0.08%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("E12").Value = "  +1.09%  "

$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.02"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").Value = "2.828.84"
$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("E16").Value = "  -1.26%  "

$ws.Range("D17").Value = "60.478.36"

$ws.Range("D18").Value = "2.402.72"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.06"
$ws.Range("E19").Value = "  +7.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.62"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.30"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("E22").Value = "  +0.70%  "

$ws.Range("E23").Value = "  -0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.81"
$ws.Range("E25").Value = "  -2.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.89"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "562.07"
$ws.Range("E27").Value = "  -3.57%  "

$ws.Range("E28").Value = "  -4.76%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").Value = "0.0₃0935"
$ws.Range("E30").Value = "  +1.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.06"

$ws.Range("E32").Value = "  -1.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.81"
$ws.Range("E33").Value = "  -2.02%  "

$ws.Range("E34").Value = "  -1.24%  "

$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("E36").Value = "  +3.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.19"
$ws.Range("E37").Value = "  +1.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.370"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("E39").Value = "  -1.03%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("E41").Value = "  -0.30%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.72"
$ws.Range("E43").Value = "  +1.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.66"
$ws.Range("E44").Value = "  -0.76%  "

$ws.Range("E45").Value = "  +5.68%  "

$ws.Range("D46").Value = "0.0₆0279"
$ws.Range("E46").Value = "  -2.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.61"
$ws.Range("E47").Value = "  +0.53%  "

$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("E49").Value = "  -0.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.27"
$ws.Range("E51").Value = "  -1.38%  "
